# issue #5: stock data from json to db
#
# The "股票" (stock) worksheet (5th sheet) is extended with three extra
# data columns (category / source_file / index) so each stock record
# matches the richer schema now produced by the json->db pipeline:
#
#   old layout (A:K) -> new layout (A:N)
#   A index   B name  C owner  D quantity  E face_value  F currency
#   G total   H property_category  [NEW] I category  (was I) J date
#   (was J) K legislator_name   (was K) L legislator_id
#   [NEW] M source_file   [NEW] N index
#
# i.e. a new "category" column is inserted right after
# "property_category" (pushing date / legislator_name / legislator_id one
# column to the right), and two brand-new trailing columns "source_file"
# and "index" are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

$lastRow = 11

# Insert a new column I ("category"), shifting the former I:K
# (date / legislator_name / legislator_id) one column right to J:L.
$ws.Columns.Item(9).Insert()

# --- header row -------------------------------------------------------
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- data rows ----------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value  = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp8fef1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
